$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12/13 swap: TRON moves up to row 12, Polkadot moves down to row 13 ---
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.114"
$ws.Range("E12").Value = "  +7.05%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "7.61"
$ws.Range("E13").Value = "  -0.36%  "

# --- Row 41/42/43 rotation: RenderToken->41, NEARProtocol->42, VeChain->43 ---
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "3.89"
$ws.Range("E41").Value = "  +1.80%  "

$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").Value = "3.41"
$ws.Range("E42").Value = "  -0.44%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0310"
$ws.Range("E43").Value = "  +0.89%  "

# --- Remaining price/volume updates ---
$ws.Range("D2").Value = "42.848.66"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.582.36"
$ws.Range("E3").Value = "  +2.58%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "302.87"
$ws.Range("E5").Value = "  +2.20%  "
$ws.Range("D6").Value = "96.87"
$ws.Range("E6").Value = "  +6.28%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.549"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("D11").Value = "0.0809"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D14").Value = "2.603.28"
$ws.Range("E14").Value = "  +3.20%  "
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "14.37"
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("D17").Value = "42.912.36"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("E18").Value = "  +5.00%  "
$ws.Range("E19").Value = "  +2.99%  "
$ws.Range("D20").Value = "6.65"
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").Value = "71.91"
$ws.Range("E21").Value = "  -1.21%  "
$ws.Range("D22").Value = "254.69"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("E23").Value = "  +2.55%  "
$ws.Range("E24").Value = "  -2.14%  "
$ws.Range("D25").Value = "28.66"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'10.20"
$ws.Range("E27").Value = "  +2.72%  "
$ws.Range("D28").Value = "39.12"
$ws.Range("E28").Value = "  +6.51%  "
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").Value = "6.03"
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("D31").Value = "155.31"
$ws.Range("E31").Value = "  +2.54%  "
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").Value = "2.75"
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("E34").Value = "  +2.05%  "
$ws.Range("D35").Value = "3.36"
$ws.Range("E35").Value = "  -3.01%  "
$ws.Range("D36").Value = "18.32"
$ws.Range("E36").Value = "  +12.10%  "
$ws.Range("D37").Value = "0.114"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("D39").Value = "23.32"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D44").Value = "2.072.92"
$ws.Range("E44").Value = "  +2.61%  "
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  +4.34%  "
$ws.Range("D47").Value = "85.14"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").Value = "75.95"
$ws.Range("E48").Value = "  +11.28%  "
$ws.Range("D49").Value = "2.836.02"
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("D50").Value = "106.03"
$ws.Range("E50").Value = "  +3.12%  "
$ws.Range("E51").Value = "  +2.69%  "
